# Update NATMI LR-pairs TPM output values (Lcn2-Lrp2)
# This applies the new TPM-derived values for rows 2-7 as per the updated
# CellPhoneDB/NATMI computation pipeline ("update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.357176
$ws.Range("H2").Value = 0.714352
$ws.Range("I2").Value = 0.005997202785179311
$ws.Range("J2").Value = 0.004010340569787885
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.032708
$ws.Range("N2").Value = 0.065416
$ws.Range("Q2").Value = 0.011682512608
$ws.Range("R2").Value = 0.046730050432
$ws.Range("S2").Value = 0.005997202785179311
$ws.Range("T2").Value = 0.004010340569787885
$ws.Range("H3").Value = 9.109331000000001
$ws.Range("I3").Value = 0.05098374027026849
$ws.Range("J3").Value = 0.05113938180746529
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.032708
$ws.Range("N3").Value = 0.065416
$ws.Range("Q3").Value = 0.09931599944933335
$ws.Range("R3").Value = 0.5958959966960001
$ws.Range("S3").Value = 0.05098374027026849
$ws.Range("T3").Value = 0.05113938180746529
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.166435666666667
$ws.Range("H4").Value = 9.499307000000002
$ws.Range("I4").Value = 0.05316638519728214
$ws.Range("J4").Value = 0.0533286898433406
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.032708
$ws.Range("N4").Value = 0.065416
$ws.Range("Q4").Value = 0.1035677777853334
$ws.Range("R4").Value = 0.6214066667120002
$ws.Range("S4").Value = 0.05316638519728214
$ws.Range("T4").Value = 0.0533286898433406
$ws.Range("I5").Value = 0.003133228500602422
$ws.Range("J5").Value = 0.002095195680465199
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.032708
$ws.Range("N5").Value = 0.065416
$ws.Range("Q5").Value = 0.006103509048
$ws.Range("R5").Value = 0.024414036192
$ws.Range("S5").Value = 0.003133228500602422
$ws.Range("T5").Value = 0.002095195680465199
$ws.Range("G6").Value = 52.102181
$ws.Range("H6").Value = 156.306543
$ws.Range("I6").Value = 0.8748273820388733
$ws.Range("J6").Value = 0.8774980271857494
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.032708
$ws.Range("N6").Value = 0.065416
$ws.Range("Q6").Value = 1.704158136148
$ws.Range("R6").Value = 10.224948816888
$ws.Range("S6").Value = 0.8748273820388733
$ws.Range("T6").Value = 0.8774980271857494
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7082566666666666
$ws.Range("H7").Value = 2.12477
$ws.Range("I7").Value = 0.01189206120779433
$ws.Range("J7").Value = 0.01192836491319154
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.032708
$ws.Range("N7").Value = 0.065416
$ws.Range("Q7").Value = 0.09931599944933335
$ws.Range("R7").Value = 0.13899395432
$ws.Range("S7").Value = 0.01189206120779433
$ws.Range("T7").Value = 0.01192836491319154

Write-Output "Updated Lcn2-Lrp2 NATMI TPM values."
